$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

# Force text format so values (dates, numbers-as-strings) are kept as plain text,
# matching the original workbook where every cell is stored as a shared string.
$ws.Range("A2:M3").NumberFormat = "@"

# Row 2 - update existing student record
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "2025-03-12"
$ws.Range("C2").Value = "satyam"
$ws.Range("D2").Value = "babu ji"
$ws.Range("E2").Value = "ramkrishnanagar"
$ws.Range("F2").Value = "7250585057"
$ws.Range("G2").Value = "06:00-10:00, 22:00-06:00"
$ws.Range("H2").Value = "1"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "350.00"
$ws.Range("K2").Value = "10.00"
$ws.Range("L2").Value = "1"
$ws.Range("M2").Value = "2025-04-12"

# Row 3 - update existing student record
$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "2025-02-11"
$ws.Range("C3").Value = "test"
$ws.Range("D3").Value = "testfather"
$ws.Range("E3").Value = "ramkrishnanagar"
$ws.Range("F3").Value = "7250585058"
$ws.Range("G3").Value = "10:00-14:00"
$ws.Range("H3").Value = "1"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "350.00"
$ws.Range("K3").Value = "11.00"
$ws.Range("L3").Value = "3"
$ws.Range("M3").Value = "2025-03-06"

# Restore default cell formatting so only the values change, matching the
# original workbook (no explicit per-cell style indices).
$ws.Range("A2:M3").ClearFormats()
